$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collected Minutiae")
$finalWs = $wb.Worksheets.Item("Final Template")

# Fix X minutiae value for minutiae points 7, 13 and 14 (approximate guess
# using the provided min image). These values previously overflowed
# DEC2HEX(x,2) (>255), producing #NUM! errors that rippled through the
# downstream hex-encoding/concatenation formulas.
$ws.Range("C8").Value = 251.5
$ws.Range("C14").Value = 253.5
$ws.Range("C15").Value = 250.5

# Nudge the long-range CONCATENATE formulas so they pick up the new
# dependent values and recompute from #NUM! to their real text results.
$ws.Range("P10").Formula = $ws.Range("P10").Formula
$ws.Range("Q2").Formula = $ws.Range("Q2").Formula
$finalWs.Range("A2").Formula = $finalWs.Range("A2").Formula

# Reflect the cell selection left by the author while making the edit
$ws.Activate()
$ws.Range("C15").Select()
